$d = $word.ActiveDocument

# Locate the five target paragraphs within the "Major Research Projects" section
# by matching on their (pre-edit) text content, then update each one's text
# directly. Using paragraph indices (rather than repeated Find/Replace calls)
# avoids ambiguity caused by the new text of one paragraph matching the search
# text intended for another paragraph (two of the paragraphs swap places).

$count = $d.Paragraphs.Count
$idxLinking = -1
$idxMeritTitle = -1
$idxMeritDesc = -1
$idxIllinoisTitle = -1
$idxIllinoisDesc = -1

for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "Linking Community Colleges and High Schools in Texas (Fall 2022 - Present)*") {
        $idxLinking = $i
    }
    elseif ($t -like "Merit-Based Financial Aid (Fall 2021 - Present)*") {
        $idxMeritTitle = $i
    }
    elseif ($t -like "Created a new dataset of all state-run merit-based aid programs*") {
        $idxMeritDesc = $i
    }
    elseif ($t -like "Exploring Equitable Higher Education Funding Models in Illinois (Summer 2022)*") {
        $idxIllinoisTitle = $i
    }
    elseif ($t -like "Collected data on Illinois appropriations for higher education*") {
        $idxIllinoisDesc = $i
    }
}

# 1. Update the "Linking Community Colleges..." project date range
$d.Paragraphs.Item($idxLinking).Range.Text = `
    "Linking Community Colleges and High Schools in Texas (Fall 2022 - Spring 2023)"

# 2 & 4. Swap the "Merit-Based Financial Aid" title and the "Exploring Equitable..." title
$d.Paragraphs.Item($idxMeritTitle).Range.Text = `
    "Exploring Equitable Higher Education Funding Models in Illinois (Summer 2022)"
$d.Paragraphs.Item($idxIllinoisTitle).Range.Text = `
    "Merit-Based Financial Aid (Fall 2021 - Present)"

# 3 & 5. Swap the two project descriptions
$d.Paragraphs.Item($idxMeritDesc).Range.Text = `
    "Collected data on Illinois appropriations for higher education, labor market projections, census data, and institutional behavior, then drafted data visualizations using R (tables, graphs, and maps) for a report targeted at state policy makers."
$d.Paragraphs.Item($idxIllinoisDesc).Range.Text = `
    "Created a new dataset of all state-run merit-based aid programs in the United States including details on dates active, award amount, merit criteria, and funding rules, leading to a new typology being proposed (see Hu et al., 2023)."
